# Generate Report for Handoff
# This script updates the localization-status workbook to reflect that the
# b578155d-... and d5e47455-... files are now "Ready for handoff" (instead of
# "Handed back: in sync with en-US"), refreshes their handoff timestamps, and
# records an "Error Detail" message noting the handback file version is stale.

$wb = $excel.ActiveWorkbook

$readyStatus = "Ready for handoff"

$b578Msg = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e06180fe166a672c46c4ceb1601d8efa8891fcee/e2e/b578155d-b0bc-4d93-9646-485830bba23d.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/cd2af7a8434389ce1937c619d411a356f26f95f5/e2e/b578155d-b0bc-4d93-9646-485830bba23d.md."
$d5e4Msg = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e06180fe166a672c46c4ceb1601d8efa8891fcee/e2e/d5e47455-8918-4c58-975d-52ce5e2b11db.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/cd2af7a8434389ce1937c619d411a356f26f95f5/e2e/d5e47455-8918-4c58-975d-52ce5e2b11db.md."

# ---------------------------------------------------------------------------
# Overview sheet: rows 4 (b578155d) and 5 (d5e47455)
#   zh-cn (E) / de-de (F) status columns -> "Ready for handoff"
#   Latest HO Xliff Generate Date (G) -> new timestamp
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("E4").Value = $readyStatus
$wsOverview.Range("F4").Value = $readyStatus
$wsOverview.Range("G4").Value = "2016-09-02 02:32:11"

$wsOverview.Range("E5").Value = $readyStatus
$wsOverview.Range("F5").Value = $readyStatus
$wsOverview.Range("G5").Value = "2016-09-02 02:32:11"

# ---------------------------------------------------------------------------
# zh-cn sheet: rows 4 (b578155d) and 5 (d5e47455)
#   Status (C) -> "Ready for handoff"
#   Latest Handoff Datetime (H) -> new timestamp
#   Error Detail (P) -> stale handback version message
#   Error Detail column width -> 40
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("C4").Value = $readyStatus
$wsZhCn.Range("H4").Value = "2016-09-02 02:32:00"
$wsZhCn.Range("P4").Value = $b578Msg

$wsZhCn.Range("C5").Value = $readyStatus
$wsZhCn.Range("H5").Value = "2016-09-02 02:32:00"
$wsZhCn.Range("P5").Value = $d5e4Msg

$wsZhCn.Columns.Item(16).ColumnWidth = $wsZhCn.Columns.Item(1).ColumnWidth

# ---------------------------------------------------------------------------
# de-de sheet: rows 4 (b578155d) and 5 (d5e47455)
#   Status (C) -> "Ready for handoff"
#   Latest Handoff Datetime (H) -> new timestamp
#   Error Detail (P) -> stale handback version message
#   Error Detail column width -> 40
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("C4").Value = $readyStatus
$wsDeDe.Range("H4").Value = "2016-09-02 02:32:11"
$wsDeDe.Range("P4").Value = $b578Msg

$wsDeDe.Range("C5").Value = $readyStatus
$wsDeDe.Range("H5").Value = "2016-09-02 02:32:11"
$wsDeDe.Range("P5").Value = $d5e4Msg

$wsDeDe.Columns.Item(16).ColumnWidth = $wsDeDe.Columns.Item(1).ColumnWidth

Write-Host "Report regenerated for handoff."
